# Apply 2024-11-02 violent crime data update (autogenerated from verified diff mapping)
$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6824
$ws.Range("K3").Value = 7053
$ws.Range("H4").Value = 1742
$ws.Range("K4").Value = 1453
$ws.Range("K6").Value = 7720
$ws.Range("H7").Value = 26055
$ws.Range("K7").Value = 23550

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K4").Value = 85
$ws.Range("K6").Value = 168
$ws.Range("K7").Value = 716
$ws.Range("K8").Value = 1542
$ws.Range("K11").Value = 440
$ws.Range("K12").Value = 42
$ws.Range("K14").Value = 117
$ws.Range("K18").Value = 156
$ws.Range("K19").Value = 692
$ws.Range("K20").Value = 570
$ws.Range("K23").Value = 233
$ws.Range("K25").Value = 110
$ws.Range("K27").Value = 222
$ws.Range("K29").Value = 1279
$ws.Range("K30").Value = 89
$ws.Range("K31").Value = 261
$ws.Range("K33").Value = 1015
$ws.Range("K37").Value = 797
$ws.Range("K42").Value = 869
$ws.Range("K43").Value = 189
$ws.Range("K45").Value = 34
$ws.Range("K48").Value = 303
$ws.Range("K49").Value = 131
$ws.Range("K50").Value = 110
$ws.Range("K51").Value = 294
$ws.Range("K52").Value = 623
$ws.Range("K53").Value = 300
$ws.Range("K54").Value = 460
$ws.Range("K56").Value = 25
$ws.Range("H63").Value = 294
$ws.Range("K63").Value = 63
$ws.Range("K64").Value = 146
$ws.Range("K65").Value = 544
$ws.Range("K67").Value = 913
$ws.Range("K70").Value = 41
$ws.Range("K71").Value = 73
$ws.Range("K77").Value = 161
$ws.Range("K78").Value = 274
$ws.Range("K82").Value = 28
$ws.Range("K84").Value = 190
$ws.Range("K85").Value = 1088
$ws.Range("K88").Value = 256
$ws.Range("K89").Value = 354
$ws.Range("K92").Value = 88
$ws.Range("K94").Value = 315
$ws.Range("K96").Value = 253
$ws.Range("K98").Value = 119
$ws.Range("K100").Value = 42
$ws.Range("H101").Value = 26055
$ws.Range("K101").Value = 23550

# Sheet 3: Bridgeport
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 117

# Sheet 4: West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 52
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 253

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 239
$ws.Range("K3").Value = 230
$ws.Range("K7").Value = 716

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 112
$ws.Range("K6").Value = 149
$ws.Range("K7").Value = 440

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 99
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 354

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 380
$ws.Range("K7").Value = 1088

# Sheet 9: Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 175
$ws.Range("K6").Value = 227
$ws.Range("K7").Value = 623

# Sheet 11: Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 76
$ws.Range("K7").Value = 300

# Sheet 12: Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 431
$ws.Range("K3").Value = 469
$ws.Range("K4").Value = 88
$ws.Range("K7").Value = 1542

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 258
$ws.Range("K6").Value = 321
$ws.Range("K7").Value = 1015

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K4").Value = 37
$ws.Range("K7").Value = 797

# Sheet 17: New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 200
$ws.Range("K7").Value = 544

# Sheet 19: Fuller Park
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 89

# Sheet 20: Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 85
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 261

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 248
$ws.Range("K7").Value = 913

# Sheet 22: South Deering
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 78
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 190

# Sheet 23: Lincoln Park
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 131

# Sheet 24: Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 249
$ws.Range("K7").Value = 460

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 362
$ws.Range("K3").Value = 454
$ws.Range("K6").Value = 372
$ws.Range("K7").Value = 1279

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 46
$ws.Range("K6").Value = 143
$ws.Range("K7").Value = 303

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 204
$ws.Range("K6").Value = 228
$ws.Range("K7").Value = 692

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 45
$ws.Range("K7").Value = 168

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 235
$ws.Range("K3").Value = 263
$ws.Range("K7").Value = 869

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 82
$ws.Range("K7").Value = 274

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 233

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 146

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 197
$ws.Range("K3").Value = 184
$ws.Range("K6").Value = 154
$ws.Range("K7").Value = 570

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 156

# Sheet 49: Wrigleyville
$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 42

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 65
$ws.Range("K6").Value = 143
$ws.Range("K7").Value = 315

# Sheet 52: East Side
$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 110

# Sheet 55: Wicker Park
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 119

# Sheet 56: Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 110

# Sheet 66: West Elsdon
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 88

# Sheet 67: O'Hare
$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 41

# Sheet 68: United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 79
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 256

# Sheet 71: Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 222

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 294

# Sheet 79: Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 37
$ws.Range("K7").Value = 189

# Sheet 81: Oakland
$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 73

# Sheet 83: Sheffield & DePaul
$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K5").Value = 16
$ws.Range("K6").Value = 28

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 68
$ws.Range("K7").Value = 161

# Sheet 85: Jackson Park
$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 34

# Sheet 86: Magnificent Mile
$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 25

# Sheet 90: Archer Heights
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 85

# Sheet 91: Beverly
$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 42
